$d = $word.ActiveDocument

# --- Styles cleanup: Normal paragraph style no longer derives from the
# custom "DocDefaults" style, and that style is removed outright. ---
$docDefaults = $d.Styles("DocDefaults")
$docDefaults.Delete()
$normal = $d.Styles("Normal")
$normal.BaseStyle = $null

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# Title text
Replace-Text "Produce a table from margins" "margins"

# Row / header labels (translate to Chinese)
Replace-Text "foreign" "国籍"
Replace-Text "Domestic" "国内"
Replace-Text "Foreign" "国外"
Replace-Text "rep78" "维修记录78"

# Domestic row numeric values
Replace-Text "22.17955" "4.830705"
Replace-Text ".577852" ".1115355"
Replace-Text "38.38" "43.31"
Replace-Text "21.02444" "4.607748"
Replace-Text "23.33466" "5.053661"

# Foreign row numeric values
Replace-Text "19.25628" "5.471863"
Replace-Text "1.019649" ".19681"
Replace-Text "18.89" "27.80"
Replace-Text "17.21803" "5.078446"
Replace-Text "21.29453" "5.865281"

# rep78 = 1 row
Replace-Text "20.53792" "4.9489"
Replace-Text "2.430089" ".4690495"
Replace-Text "8.45" "10.55"
Replace-Text "15.68024" "4.011283"
Replace-Text "25.3956" "5.886517"

# rep78 = 2 row
Replace-Text "20.25927" "5.095182"
Replace-Text "1.248904" ".2410603"
Replace-Text "16.22" "21.14"
Replace-Text "17.76275" "4.613309"
Replace-Text "22.7558" "5.577054"

# rep78 = 3 row
Replace-Text "20.5155" "5.094224"
Replace-Text ".6630665" ".1279834"
Replace-Text "30.94" "39.80"
Replace-Text "19.19005" "4.838389"
Replace-Text "21.84095" "5.35006"

# rep78 = 4 row
Replace-Text "21.21928" "5.003851"
Replace-Text ".8307747" ".160354"
Replace-Text "25.54" "31.21"
Replace-Text "19.55858" "4.683308"
Replace-Text "22.87997" "5.324394"

# rep78 = 5 row
Replace-Text "24.40345" "4.838878"
Replace-Text "1.185693" ".2288594"
Replace-Text "20.58" "21.14"
Replace-Text "22.03329" "4.381394"
Replace-Text "26.77362" "5.296361"
